$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.515.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.32%  "

$ws.Range("D3").Value = "'1.813.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.37%  "

$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").Value = "'225.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.80%  "

$ws.Range("E6").Value = "  +3.01%  "

$ws.Range("E7").Value = "  -0.22%  "

$ws.Range("D8").Value = "'38.06"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.19%  "

$ws.Range("E9").Value = "  -3.71%  "

$ws.Range("D10").Value = "'0.0677"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.53%  "

$ws.Range("D11").Value = "'0.0974"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.98%  "

$ws.Range("D12").Value = "'2.076.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.30%  "

$ws.Range("D13").Value = "'11.24"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.76%  "

$ws.Range("D14").Value = "'1.835.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.55%  "

$ws.Range("E15").Value = "  -1.76%  "

$ws.Range("D16").Value = "'34.499.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.23%  "

$ws.Range("D17").Value = "'4.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.92%  "

$ws.Range("D18").Value = "'68.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.90%  "

$ws.Range("D19").Value = "'243.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.03%  "

$ws.Range("D20").Value = "'0.0₃0775"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.86%  "

$ws.Range("D21").Value = "'11.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.87%  "

$ws.Range("E22").Value = "  -0.22%  "

$ws.Range("E23").Value = "  -1.41%  "

$ws.Range("E24").Value = "  +3.84%  "

$ws.Range("D25").Value = "'170.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.67%  "

$ws.Range("D26").Value = "'7.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.58%  "

$ws.Range("D27").Value = "'17.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.35%  "

$ws.Range("E28").Value = "  +1.09%  "

$ws.Range("E29").Value = "  -0.24%  "

$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "'3.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.56%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.03%  "

$ws.Range("D32").Value = "'0.0520"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.37%  "

$ws.Range("D33").Value = "'3.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.02%  "

$ws.Range("E34").Value = "  +0.12%  "

$ws.Range("D35").Value = "'1.366.13"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.40%  "

$ws.Range("D36").Value = "'0.647"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.01%  "

$ws.Range("D37").Value = "'1.06"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.36%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0187"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.01%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'2.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.79%  "

$ws.Range("E40").Value = "  -2.68%  "

$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").Value = "'2.44"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.16%  "

$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'0.955"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.72%  "

$ws.Range("D43").Value = "'81.76"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.29%  "

$ws.Range("E44").Value = "  -1.26%  "

$ws.Range("D45").Value = "'13.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.06%  "

$ws.Range("D46").Value = "'0.0509"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.68%  "

$ws.Range("D47").Value = "'1.976.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.32%  "

$ws.Range("D48").Value = "'5.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.08%  "

$ws.Range("E49").Value = "  -0.24%  "

$ws.Range("D50").Value = "'102.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.89%  "

$ws.Range("E51").Value = "  -5.35%  "
